# Auto-generated Excel COM-interop script
# Applies updated pricing / profit values to the 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as produced by the scheduled market-data refresh run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2372.2
$ws.Range("I18").Value = 2540.7058
$ws.Range("K18").Value = 2540.7058
$ws.Range("M18").Value = -2256.7058

$ws.Range("H39").Value = 1841.7646
$ws.Range("I39").Value = 345.45456
$ws.Range("K39").Value = 1036.36368
$ws.Range("M39").Value = -740.3636799999999

$ws.Range("H40").Value = 5312.5
$ws.Range("I40").Value = 4164.8887
$ws.Range("J40").Value = 6460.1113
$ws.Range("K40").Value = 4164.8887
$ws.Range("L40").Value = 6460.1113
$ws.Range("M40").Value = -3989.8887
$ws.Range("N40").Value = -6810.1113

$ws.Range("H74").Value = 6527.5
$ws.Range("I74").Value = 4507.4614
$ws.Range("J74").Value = 9445.333000000001
$ws.Range("K74").Value = 4507.4614
$ws.Range("L74").Value = 9445.333000000001
$ws.Range("M74").Value = -3571.4614
$ws.Range("N74").Value = -11317.333

$ws.Range("H77").Value = 6527.5
$ws.Range("I77").Value = 4507.4614
$ws.Range("J77").Value = 9445.333000000001
$ws.Range("K77").Value = 22537.307
$ws.Range("L77").Value = 47226.665
$ws.Range("M77").Value = -17857.307
$ws.Range("N77").Value = -56586.665

$ws.Range("H98").Value = 1088.7805
$ws.Range("I98").Value = 1007.51514
$ws.Range("K98").Value = 1007.51514
$ws.Range("M98").Value = 490.48486

$ws.Range("H116").Value = 4000
$ws.Range("I116").Value = 4000
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 4000
$ws.Range("L116").Value = 4000
$ws.Range("M116").Value = -558
$ws.Range("N116").Value = -10884

$ws.Range("H122").Value = 1088.7805
$ws.Range("I122").Value = 1007.51514
$ws.Range("K122").Value = 3022.54542
$ws.Range("M122").Value = -572.5454199999999

$ws.Range("H125").Value = 2768
$ws.Range("I125").Value = 2248.25
$ws.Range("K125").Value = 20234.25
$ws.Range("M125").Value = -17774.25

$ws.Range("H132").Value = 2199.0557
$ws.Range("I132").Value = 1242.75
$ws.Range("K132").Value = 3728.25
$ws.Range("M132").Value = -1198.25

$ws.Range("H141").Value = 1693.0588
$ws.Range("I141").Value = 1590.6154
$ws.Range("K141").Value = 4771.8462
$ws.Range("M141").Value = 408.1538

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3299240.5
$ws.Range("I32").Value = 3658066.8
$ws.Range("K32").Value = 3658066.8
$ws.Range("M32").Value = -3657779.8

$ws.Range("H61").Value = 6864.1035
$ws.Range("I61").Value = 6556.087
$ws.Range("K61").Value = 6556.087
$ws.Range("M61").Value = -6344.087

$ws.Range("H74").Value = 457382.3
$ws.Range("I74").Value = 589462.7
$ws.Range("K74").Value = 589462.7
$ws.Range("M74").Value = -588588.7

$ws.Range("H77").Value = 457382.3
$ws.Range("I77").Value = 589462.7
$ws.Range("K77").Value = 2947313.5
$ws.Range("M77").Value = -2942945.5

$ws.Range("H97").Value = 952999.6
$ws.Range("I97").Value = 1280988.9
$ws.Range("J97").Value = 1830.8
$ws.Range("K97").Value = 1280988.9
$ws.Range("L97").Value = 1830.8
$ws.Range("M97").Value = -1280492.9
$ws.Range("N97").Value = -2822.8

$ws.Range("H122").Value = 2373
$ws.Range("I122").Value = 2301.8333
$ws.Range("K122").Value = 6905.499899999999
$ws.Range("M122").Value = -4455.499899999999

$ws.Range("H125").Value = 50715
$ws.Range("J125").Value = 50715
$ws.Range("L125").Value = 50715
$ws.Range("M125").Value = -60555

$ws.Range("H136").Value = 6864.1035
$ws.Range("I136").Value = 6556.087
$ws.Range("K136").Value = 19668.261
$ws.Range("M136").Value = -17118.261

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 1145
$ws.Range("I44").Value = 1145
$ws.Range("K44").Value = 1145
$ws.Range("M44").Value = -648

$ws.Range("H134").Value = 6276.1
$ws.Range("I134").Value = 2429.6897
$ws.Range("K134").Value = 7289.0691
$ws.Range("M134").Value = -4754.0691

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3232.0715
$ws.Range("I22").Value = 544.6
$ws.Range("K22").Value = 544.6
$ws.Range("M22").Value = -194.6

$ws.Range("H31").Value = 32262290
$ws.Range("I31").Value = 58825964
$ws.Range("K31").Value = 58825964
$ws.Range("M31").Value = -58825669

$ws.Range("H34").Value = 32262290
$ws.Range("I34").Value = 58825964
$ws.Range("K34").Value = 58825964
$ws.Range("M34").Value = -58825762

$ws.Range("H122").Value = 84994.5
$ws.Range("I122").Value = 251189.5
$ws.Range("K122").Value = 753568.5
$ws.Range("M122").Value = -751118.5

$ws.Range("H132").Value = 5170.8857
$ws.Range("I132").Value = 3382.8333
$ws.Range("K132").Value = 10148.4999
$ws.Range("M132").Value = -7618.499899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4368009.5
$ws.Range("I4").Value = 3966770.5
$ws.Range("J4").Value = 5055848
$ws.Range("K4").Value = 11900311.5
$ws.Range("L4").Value = 15167544
$ws.Range("M4").Value = -11900199.5
$ws.Range("N4").Value = -15167768

$ws.Range("H5").Value = 951.3333
$ws.Range("I5").Value = 877
$ws.Range("J5").Value = 1100
$ws.Range("K5").Value = 2631
$ws.Range("L5").Value = 3300
$ws.Range("M5").Value = -2519
$ws.Range("N5").Value = -3524

$ws.Range("H62").Value = 2110.6667

$ws.Range("H65").Value = 2110.6667

$ws.Range("H69").Value = 3994.5
$ws.Range("I69").Value = 1099.5
$ws.Range("K69").Value = 3298.5
$ws.Range("M69").Value = -2487.5

$ws.Range("H72").Value = 3994.5
$ws.Range("I72").Value = 1099.5
$ws.Range("K72").Value = 9895.5
$ws.Range("M72").Value = -5839.5

$ws.Range("H80").Value = 3354.8333
$ws.Range("J80").Value = 3354.8333
$ws.Range("L80").Value = 10064.4999
$ws.Range("N80").Value = -11936.4999

$ws.Range("H83").Value = 3354.8333
$ws.Range("J83").Value = 3354.8333
$ws.Range("L83").Value = 30193.4997
$ws.Range("N83").Value = -39553.4997

$ws.Range("H103").Value = 2082.6
$ws.Range("I103").Value = 138.33333
$ws.Range("K103").Value = 414.99999
$ws.Range("M103").Value = 464.00001

$ws.Range("H135").Value = 951.3333
$ws.Range("I135").Value = 877
$ws.Range("J135").Value = 1100
$ws.Range("K135").Value = 7893
$ws.Range("L135").Value = 9900
$ws.Range("M135").Value = -5358
$ws.Range("N135").Value = -14970

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 143486.28
$ws.Range("I12").Value = 126488.75
$ws.Range("K12").Value = 126488.75
$ws.Range("M12").Value = -126348.75

$ws.Range("H107").Value = 1535.75
$ws.Range("J107").Value = 1470.1666
$ws.Range("L107").Value = 1470.1666
$ws.Range("N107").Value = -5310.1666

$ws.Range("H113").Value = 32119.4
$ws.Range("J113").Value = 2599.3333
$ws.Range("L113").Value = 2599.3333
$ws.Range("N113").Value = -6939.3333

$ws.Range("H132").Value = 4422.325
$ws.Range("I132").Value = 2251.6553
$ws.Range("J132").Value = 10145
$ws.Range("K132").Value = 6754.965899999999
$ws.Range("L132").Value = 30435
$ws.Range("M132").Value = -4224.965899999999
$ws.Range("N132").Value = -35495

$ws.Range("H140").Value = 120000
$ws.Range("J140").Value = 120000
$ws.Range("L140").Value = 120000
$ws.Range("N140").Value = -130360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4652.316
$ws.Range("J46").Value = 4916.7715
$ws.Range("L46").Value = 4916.7715
$ws.Range("N46").Value = -5292.7715

$ws.Range("H55").Value = 366.9375
$ws.Range("I55").Value = 161.55556
$ws.Range("J55").Value = 631
$ws.Range("K55").Value = 161.55556
$ws.Range("L55").Value = 631
$ws.Range("M55").Value = 11.44443999999999
$ws.Range("N55").Value = -977

$ws.Range("H136").Value = 4243.2
$ws.Range("I136").Value = 3433.1765
$ws.Range("K136").Value = 10299.5295
$ws.Range("M136").Value = -7749.529500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 170002.33
$ws.Range("I8").Value = 500015
$ws.Range("K8").Value = 500015
$ws.Range("M8").Value = -499875

$ws.Range("H93").Value = 78000
$ws.Range("J93").Value = 78000
$ws.Range("L93").Value = 78000
$ws.Range("N93").Value = -82992

$ws.Range("H132").Value = 3118.2964
$ws.Range("I132").Value = 2310.4634
$ws.Range("K132").Value = 6931.3902
$ws.Range("M132").Value = -4401.3902

$ws.Range("H136").Value = 3155.3333
$ws.Range("I136").Value = 937.05884
$ws.Range("K136").Value = 2811.17652
$ws.Range("M136").Value = -261.17652
